# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# Map of row -> new F value, shared between the two affected sheets.
# (Row numbers differ slightly between sheets because "全部类型" has two
# extra rows before the last block of changed rows.)
$updatesSheet1 = @{
    4  = 1785
    5  = 35
    7  = 664
    18 = 108
    19 = 4968
    23 = 2244
    24 = 71
    25 = 24
    26 = 2094
}

$updatesSheet4 = @{
    4  = 1785
    5  = 35
    7  = 664
    18 = 108
    19 = 4968
    25 = 2244
    26 = 71
    27 = 24
    28 = 2094
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updatesSheet1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updatesSheet1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesSheet4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updatesSheet4[$row]
}
